$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update ID_SOCIO values (column A) for rows 2-4
$ws.Range("A2").Value = 103
$ws.Range("A3").Value = 105
$ws.Range("A4").Value = 101

# Update NOMBRE / APELLIDOS (columns B and C) for rows 2-4
$ws.Range("B2").Value = "Eduardo"
$ws.Range("C2").Value = "Educado"

$ws.Range("B3").Value = "Lola"
$ws.Range("C3").Value = "Sol"

$ws.Range("B4").Value = "Gustavo"
$ws.Range("C4").Value = "Rana"
